$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row (row 24) produced by the 2025-12-18 run.
$row = 24

# Column A holds a literal date-like string (matches existing rows, which
# store "MM/DD/YYYY" as plain text, not a date serial). Force text entry by
# temporarily applying a text number format, then clear the format again so
# the cell doesn't end up with leftover style/format baggage that the other
# rows in this column don't have.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "12/18/2025"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = 11475.31
$ws.Cells.Item($row, 3).Value = 0.2143824212507974
$ws.Cells.Item($row, 4).Value = 0.7856175787492026
$ws.Cells.Item($row, 5).Value = -157.74
$ws.Cells.Item($row, 6).Value = -34.29
$ws.Cells.Item($row, 7).Value = -21501.55
$ws.Cells.Item($row, 8).Value = -70.45999999999999
$ws.Cells.Item($row, 9).Value = -492.74
$ws.Cells.Item($row, 10).Value = -16.69
